# EndofsprintPresentation.pptx edit
#
#   1. Slide-number placeholder bracket text "<#>" -> "<nr.>" on the
#      slide master and on every slide layout (18 locations total).
#   2. Resized/repositioned the PCA content placeholder + pictures on
#      slide 6 ("PCA") and slide 7 ("PCA").
#
# NOTE on EMU precision: Shape.Left/Top/Width/Height are expressed in
# points, and this COM host stores the underlying value as a 32-bit
# float before converting back to EMU (truncating, not rounding), so a
# naive "EMU / 12700" assignment can land 1 EMU short of the target.
# EmuToPoints searches for the nearest representable point value whose
# float32 round-trip truncates to exactly the requested EMU count.

function EmuToPoints {
    param([long]$Emu)

    $base = [double]$Emu / 12700.0

    for ($i = 0; $i -lt 20000; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $asSingle = [float]$candidate
        $backToEmu = [int64]([double]$asSingle * 12700.0)
        if ($backToEmu -eq $Emu) {
            return $candidate
        }
    }

    return $base
}

function Set-ShapeEmuBounds {
    param($Shape, $Left, $Top, $Width, $Height)

    $Shape.Left   = EmuToPoints $Left
    $Shape.Top    = EmuToPoints $Top
    $Shape.Width  = EmuToPoints $Width
    $Shape.Height = EmuToPoints $Height
}

function Set-SlideNumberBracketText {
    param($Shapes, [string]$NewText)

    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shape = $Shapes.Item($i)
        $isSlideNumberPlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 13) {
                # ppPlaceholderSlideNumber
                $isSlideNumberPlaceholder = $true
            }
        } catch {
            $isSlideNumberPlaceholder = $false
        }

        if ($isSlideNumberPlaceholder) {
            $shape.TextFrame.TextRange.Text = $NewText
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1. "<#>" -> "<nr.>" on the slide master and every slide layout ---
$bracketText = "$([char]0x2039)nr.$([char]0x203A)"

$master = $p.SlideMaster
Set-SlideNumberBracketText $master.Shapes $bracketText

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-SlideNumberBracketText $layouts.Item($li).Shapes $bracketText
}

# --- 2. PCA slides resize/reposition ---

# Slide 6 ("PCA"): Content Placeholder 2, Picture 4, Picture 5
$slide6 = $p.Slides.Item(6)
Set-ShapeEmuBounds $slide6.Shapes.Item(2) 1104293 1331259 8946541 4195481
Set-ShapeEmuBounds $slide6.Shapes.Item(3) 646111  3204519 4852989 3134704
Set-ShapeEmuBounds $slide6.Shapes.Item(4) 6269567 3204519 4818140 3134704

# Slide 7 ("PCA"): Content Placeholder 2, Picture 6 (Picture 4 unchanged)
$slide7 = $p.Slides.Item(7)
Set-ShapeEmuBounds $slide7.Shapes.Item(2) 1104900 1459793 8946541 4195481
Set-ShapeEmuBounds $slide7.Shapes.Item(4) 1104899 2143742 8607511 4195481
